$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 new quarter columns before the existing data (E:I), pushing the
#    current E:I data (and its formatting) to J:N.
$ws.Columns("E:I").Insert()

# 2) New quarter headers (both header rows: 8 and 24)
$ws.Range("E8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("F8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("G8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("H8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("I8").Value = "فصل دوم منتهی به 1400/06"

$ws.Range("E24").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("F24").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("G24").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("H24").Value = "فصل اول منتهی به 1400/03"
$ws.Range("I24").Value = "فصل دوم منتهی به 1400/06"

# 3) New quarterly figures for the expense table (rows with all-zero series
#    keep zeros for the newly inserted quarters too)
$ws.Range("E10:I10").Value = 0
$ws.Range("E11:I11").Value = 0
$ws.Range("E12:I12").Value = 0
$ws.Range("E13:I13").Value = 0
$ws.Range("E14:I14").Value = 0

$ws.Range("E15").Value = 38791
$ws.Range("F15").Value = -176
$ws.Range("G15").Value = -38210
$ws.Range("H15").Value = 570
$ws.Range("I15").Value = 403

$ws.Range("E16").Value = 10490
$ws.Range("F16").Value = 11804
$ws.Range("G16").Value = 11855
$ws.Range("H16").Value = 21970
$ws.Range("I16").Value = 24295

$ws.Range("E17").Value = 160390
$ws.Range("F17").Value = 200318
$ws.Range("G17").Value = 290435
$ws.Range("H17").Value = 331635
$ws.Range("I17").Value = 399507

$ws.Range("E18:I18").Value = 0

$ws.Range("E19").Value = 315185
$ws.Range("F19").Value = 389940
$ws.Range("G19").Value = 563068
$ws.Range("H19").Value = 600266
$ws.Range("I19").Value = 639643

$ws.Range("E20").Value = 524856
$ws.Range("F20").Value = 601886
$ws.Range("G20").Value = 827148
$ws.Range("H20").Value = 954441
$ws.Range("I20").Value = 1063848

# 4) Headcount table: recompute the newly inserted quarters' figures, and
#    fix up the existing quarters to their updated values
$ws.Range("E26").Value = 571
$ws.Range("F26").Value = 571
$ws.Range("G26").Value = 566
$ws.Range("H26").Value = 566
$ws.Range("I26").Value = 559

$ws.Range("E27").Value = 971
$ws.Range("F27").Value = 971
$ws.Range("G27").Value = 973
$ws.Range("H27").Value = 973
$ws.Range("I27").Value = 959
